# Applies the "added support for models trained on colab, train-info updated" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 updates ---
# A3: -1 -> 1
$ws.Range("A3").Value = 1

# New metric columns H3:M3
$ws.Range("H3").Value = 0.7305
$ws.Range("I3").Value = 0.6565
$ws.Range("J3").Value = 0.7344
$ws.Range("K3").Value = 0.7697
$ws.Range("L3").Value = 0.7022
$ws.Range("M3").Value = 0.67

# N3: Env label changes from rtx5071 -> rtx5070
$ws.Range("N3").Value = "rtx5070"

# --- Row 5 updates ---
# A5: 0 -> 1
$ws.Range("A5").Value = 1

# New metric columns H5:M5
$ws.Range("H5").Value = 0.8423
$ws.Range("I5").Value = 0.8191
$ws.Range("J5").Value = 0.8442
$ws.Range("K5").Value = 0.844
$ws.Range("L5").Value = 0.8444
$ws.Range("M5").Value = 0.8222

# --- View changes ---
$ws.Application.ActiveWindow.Zoom = 130
$ws.Range("J8").Select()
